# Update the NEIGHBORHOOD_NAME column (B) on Sheet1 so adjacent ZCTAs are
# consolidated into the same, simplified neighborhood groupings.
#
# New text per ZCTA (column A holds the ZCTA, row 2..38).
$newNames = @{
    11201 = "Fort Greene/Brooklyn Hts"
    11203 = "East Flatbush"
    11204 = "Bensonhurst"
    11205 = "Fort Greene/Brooklyn Hts"
    11206 = "Williamsburg/Greenpoint"
    11207 = "East New York"
    11208 = "East New York"
    11209 = "Bay Ridge"
    11210 = "Flatbush/Midwood"
    11211 = "Williamsburg/Greenpoint"
    11212 = "Brownsville"
    11213 = "Crown Heights North"
    11214 = "Bensonhurst"
    11215 = "Park Slope"
    11216 = "Bedford Stuyvesant"
    11217 = "Park Slope"
    11218 = "Borough Park"
    11219 = "Borough Park"
    11220 = "Sunset Park"
    11221 = "Bushwick"
    11222 = "Williamsburg/Greenpoint"
    11223 = "Sheepshead Bay"
    11224 = "Coney Island"
    11225 = "Crown Heights South"
    11226 = "Flatbush/Midwood"
    11228 = "Bay Ridge"
    11229 = "Sheepshead Bay"
    11230 = "Flatbush/Midwood"
    11231 = "Park Slope"
    11232 = "Sunset Park"
    11233 = "Brownsville"
    11234 = "Canarsie"
    11235 = "Sheepshead Bay"
    11236 = "Canarsie"
    11237 = "Williamsburg/Greenpoint"
    11238 = "Crown Heights North"
    11239 = "East New York"
}

# Rows whose B cell loses its explicit (12pt) cell style and reverts to the
# workbook default ("Normal") formatting as part of this edit.
$plainRows = @(2, 4, 5, 6, 7, 8, 10, 11, 13, 14, 16, 17, 18, 19, 22, 23, 25, 26, 27, 29, 30, 31, 32, 34, 36, 37)

# Sheet row numbers in the exact order they were (re)typed by the original
# author, so freshly introduced neighborhood labels are appended to the
# shared-string table in the same sequence as in the source edit.
$writeOrder = @(
    2, 5, 4, 14, 6, 11, 22, 36, 10, 26, 29, 16, 13, 37, 25,
    3, 7, 8, 9, 12, 15, 17, 18, 19, 20, 21, 23, 24, 27, 28,
    30, 31, 32, 33, 34, 35, 38
)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($row in $writeOrder) {
    $zcta = [int]$ws.Cells.Item($row, 1).Value()
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $newNames[$zcta]
    if ($plainRows -contains $row) {
        $cell.Style = "Normal"
    }
}

$ws.Range("A2").Select()
